# "defined R environment for download-gbif script"
#
# Rename Sheet1 -> "notes on datasets" and append the new species rows
# (319-329) that were added to column A, then leave the workbook focused
# on that sheet/cell, matching the author's final selection state.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "notes on datasets"

$newSpecies = @(
    "Solanum melongena",
    "Solanum lycopersicoides",
    "Solanum pimpinellifolium",
    "Iochroma cyaneum",
    "Nicotiana attenuata",
    "Nicotiana benthamiana",
    "Nicotiana tabacum",
    "Petunia axillaris",
    "Petunia inflata",
    "Solanum chilense",
    "Coffea humblotiana"
)

$startRow = 319
for ($i = 0; $i -lt $newSpecies.Length; $i++) {
    $row = $startRow + $i
    $ws1.Cells.Item($row, 1).Value = $newSpecies[$i]
}

# Make "notes on datasets" the active sheet/tab and leave the selection on
# the last cell entered, scrolled so that row is in view.
[void]$ws1.Activate()
[void]$ws1.Range("A329").Select()
